$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells whose new text looks like a pure number need NumberFormat "@"
# applied first so they are stored as text (matching original inlineStr text cells)
# rather than being auto-converted to numeric values by Excel.

$ws.Range("D2").Value = "76.916.09"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "2.962.36"
$ws.Range("E3").Value = "  +3.48%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "200.69"
$ws.Range("E5").Value = "  +2.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "597.10"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.549"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").Value = "2.962.07"
$ws.Range("E10").Value = "  +3.41%  "
$ws.Range("E11").Value = "  +14.35%  "
$ws.Range("D13").Value = "3.511.69"
$ws.Range("E13").Value = "  +3.48%  "
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "76.821.48"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.32"
$ws.Range("E16").Value = "  +3.24%  "
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "2.953.37"
$ws.Range("E18").Value = "  +2.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.50"
$ws.Range("E19").Value = "  +8.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.70"
$ws.Range("E20").Value = "  -3.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.67"
$ws.Range("E21").Value = "  -1.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.34"
$ws.Range("E22").Value = "  +5.23%  "
$ws.Range("E23").Value = "  -3.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.84"
$ws.Range("E24").Value = "  +1.79%  "
$ws.Range("D25").Value = "3.116.26"
$ws.Range("E25").Value = "  +2.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.71"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  +4.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.34"
$ws.Range("E31").Value = "  +8.36%  "
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "499.60"
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  +2.00%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("B36").Value = "Cronos"
$ws.Range("C36").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").Value = "  +23.78%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.21"
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.397"
$ws.Range("E38").Value = "  +15.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.21"
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("E40").Value = "  +1.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.110"
$ws.Range("E41").Value = "  -5.78%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "180.88"
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.93"
$ws.Range("E44").Value = "  -1.90%  "
$ws.Range("E45").Value = "  -1.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.15"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("E47").Value = "  -2.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.592"
$ws.Range("E48").Value = "  +2.91%  "
$ws.Range("E49").Value = "  +4.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.32"
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.69"
$ws.Range("E51").Value = "  +5.66%  "
